$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.46"
$ws.Range("E2").Value = "'5.52%"
$ws.Range("D3").Value = "'31.85"
$ws.Range("E3").Value = "'7.72%"
$ws.Range("D4").Value = "'5.225"
$ws.Range("E4").Value = "'2.68%"
$ws.Range("D5").Value = "'0.07361"
$ws.Range("E5").Value = "'8.85%"
$ws.Range("D6").Value = "'7.831"
$ws.Range("E6").Value = "'6.34%"
$ws.Range("D7").Value = "'3.732"
$ws.Range("E7").Value = "'8.42%"
$ws.Range("D8").Value = "'1.504"
$ws.Range("E8").Value = "'7.15%"
$ws.Range("D9").Value = "'0.9071"
$ws.Range("E9").Value = "'-0.89%"
$ws.Range("D10").Value = "'0.01681"
$ws.Range("E10").Value = "'2,489.70%"
$ws.Range("D11").Value = "'0.1680"
$ws.Range("E11").Value = "'5.01%"
$ws.Range("D12").Value = "'0.07504"
$ws.Range("E12").Value = "'8.03%"
$ws.Range("D13").Value = "'0.07964"
$ws.Range("E13").Value = "'3.75%"
$ws.Range("E14").Value = "'0.94%"
$ws.Range("D15").Value = "'0.09914"
$ws.Range("E15").Value = "'10.14%"
$ws.Range("D16").Value = "'0.001487"
$ws.Range("E16").Value = "'-5.87%"
$ws.Range("D17").Value = "'0.04521"
$ws.Range("E17").Value = "'1.05%"
$ws.Range("D18").Value = "'0.006200"
$ws.Range("E18").Value = "'0.58%"
$ws.Range("D19").Value = "'3.470"
$ws.Range("E19").Value = "'0.55%"
$ws.Range("E20").Value = "'-0.05%"
$ws.Range("D21").Value = "'0.3338"
$ws.Range("E21").Value = "'4.35%"
$ws.Range("D22").Value = "'0.1322"
$ws.Range("E22").Value = "'1.18%"
$ws.Range("D23").Value = "'4.396"
$ws.Range("E23").Value = "'7.46%"
$ws.Range("E24").Value = "'2.38%"
$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'1.70%"
$ws.Range("D26").Value = "'0.004427"
$ws.Range("E26").Value = "'7.01%"
$ws.Range("D27").Value = "'0.0001297"
$ws.Range("E27").Value = "'8.14%"
$ws.Range("D28").Value = "'0.0001736"
$ws.Range("E28").Value = "'7.38%"
$ws.Range("D40").Value = "'0.04501"
$ws.Range("E40").Value = "'5.59%"
$ws.Range("D41").Value = "'0.007186"
$ws.Range("E41").Value = "'5.23%"
$ws.Range("D42").Value = "'0.1343"
$ws.Range("E42").Value = "'7.97%"
$ws.Range("D43").Value = "'0.002324"
$ws.Range("E43").Value = "'4.30%"
$ws.Range("D44").Value = "'0.01340"
$ws.Range("E44").Value = "'3.82%"
$ws.Range("D45").Value = "'0.00006065"
$ws.Range("E45").Value = "'7.13%"
$ws.Range("E46").Value = "'-3.45%"
$ws.Range("D47").Value = "'0.01296"
$ws.Range("E47").Value = "'-13.91%"
